# Applies the automatic update dated 2025-08-25 14:20:08
# Updates INODOROS/LAVABOS related figures for ALMEIDA CUATIN JHONATHANN CARLOS
# across the three sheets, keeping all derived totals/percentages consistent.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (row 19: MANCHENO PINO HERVIN SANTIAGO) ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("H19").Value = 586.8
$wsGrupo.Range("I19").Value = 174.6

# --- Sheet "VENTA MENSUAL" (row 19: MANCHENO PINO HERVIN SANTIAGO; row 34: TOTAL) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F19").Value = 3366.47
$wsMensual.Range("F34").Value = 16667.86

# --- Sheet "CUMPLIMIENTO MENSUAL" (row 7: INODOROS; row 8: LAVABOS; row 19: TOTAL) ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D7").Value = 586.8
$wsCumplimiento.Range("E7").Value = 1013.2
$wsCumplimiento.Range("F7").Value = 0.36675

$wsCumplimiento.Range("D8").Value = 174.6
$wsCumplimiento.Range("E8").Value = 450.4
$wsCumplimiento.Range("F8").Value = 0.27936

$wsCumplimiento.Range("D19").Value = 17186.07
$wsCumplimiento.Range("E19").Value = 14923.21107555788
$wsCumplimiento.Range("F19").Value = 0.535236835716086
